# New submission synced: 2026-02-09 17:45:32
# Sheet "JSS 3E" (Class register) gets a new response row appended, and a
# previously mis-typed Admission No cell (C11) is corrected to a real number.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("JSS 3E")

# --- Fix existing row 11: Admission No "38" was stored as text, make it numeric ---
$ws.Range("C11").Value = 38

# --- Append new submission as row 12 ---
$ws.Range("A12").Value = "2026-02-09 17:45:32"
$ws.Range("B12").Value = "Umar Fatima Ali"

# Admission No "44" is kept as text (matches source form data), so force text
# formatting before assigning the value...
$ws.Range("C12").NumberFormat = "@"
$ws.Range("C12").Value = "44"

$ws.Range("D12").Value = 10

# ...then strip the one-off "@" style back off C12 by copying the plain,
# unstyled formatting from a neighboring data cell, so the new row matches
# the rest of the table (no explicit cell style).
$ws.Range("A11").Copy()
$ws.Range("C12").PasteSpecial(-4122)  # xlPasteFormats
